$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SONGS cell in row 3 used to hold one combined list of songs for the
# "40th Year Anniversary" service. Split it into one row per song, each
# keeping the same sermon topic.
#
# Insert 3 new rows below row 3, copying row 3's formatting down so the new
# rows pick up the existing table/date-column styling (borders etc.)
# instead of plain default formatting.
$ws.Rows("3").Copy()
$ws.Rows("4:6").Insert(-4121)

# Re-apply the SONGS/SERMON TOPIC column formatting (row 3's) to the new
# rows so they match the rest of the table rather than the date column's
# formatting that Insert copied into B:C.
$ws.Range("B3:C3").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)
$ws.Range("B5:C5").PasteSpecial(-4122)
$ws.Range("B6:C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in each song individually, in the order they originally appeared in
# the combined string, except "Be Thou My Vision" (row 3) last - matching
# the order new values were entered in the source edit.
$ws.Range("B4").Value = "Tsis Muaj Koj Pab (#148)"
$ws.Range("C4").Value = "40th Year Anniversary"

$ws.Range("B5").Value = "You Raise Me Up"
$ws.Range("C5").Value = "40th Year Anniversary"

$ws.Range("B6").Value = "10,000 Reasons"
$ws.Range("C6").Value = "40th Year Anniversary"

$ws.Range("B3").Value = "Be Thou My Vision"
$ws.Range("C3").Value = "40th Year Anniversary"

# Keep the table/list object in sync with the newly added rows.
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:C6"))

# Match the saved selection state from the authored workbook.
$ws.Range("A3").Select()
